# ScanOnline_LoginData_DataDrivenTest.xlsx update (5 Oct 25):
#   Row 2 credentials changed from YuvrajScan/Admin@2029 to Admin/admin123.
#   The B2 hyperlink still points at mailto:Admin@2029 but now needs an
#   explicit display-text override since the cell text no longer matches
#   the mailto address. The saved cursor position moved from H13 to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the row 2 username value -----------------------------------
$ws.Range("A2").Value = "Admin"

# --- 2. Stash the current (non-hyperlink) formatting of column B so it can
#        be restored after the hyperlinks are rebuilt below. Adding a
#        hyperlink through the object model reformats the cell with the
#        built-in "Hyperlink" style, which this workbook does not use.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- 3. Rebuild the hyperlinks. Hyperlinks.Delete() clears every hyperlink
#        on the worksheet (there is no working single-item delete in this
#        object model), so recreate all five in their original order and
#        targets; only B2 needs a TextToDisplay override now. (Passing
#        TextToDisplay also sets the cell's text, so the real password
#        value is (re)applied to B2 afterwards, in step 4.)
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Admin@2029", "", "", "Admin@2029") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Admin@3456", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Admin@0987", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:Section@1456", "", "", "") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:GoodQa@123", "", "", "") | Out-Null

# --- 4. Set the real row 2 password value (overwrites the TextToDisplay
#        text stamped onto B2 by step 3, while keeping the display override
#        recorded on the hyperlink itself).
$ws.Range("B2").Value = "admin123"

# --- 5. Restore the original formatting on B2:B6, then clean up the
#        scratch cell used to stash it.
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("B2:B6").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ1").EntireColumn.Delete() | Out-Null
$excel.CutCopyMode = 0

# --- 6. Move the saved selection/cursor to D11 ----------------------------
$ws.Range("D11").Select() | Out-Null
